$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 matching style of G1 ("sum" header)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data column H2:H8
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 0
